$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (row 1) to the new labels used for academic advising
$ws.Range("A1").Value = "Student Name"
$ws.Range("B1").Value = "Student Mail"
$ws.Range("C1").Value = "Instructor Name "
$ws.Range("D1").Value = "Instructor Mail "

# Update the active selection to match the new view state
$ws.Range("N15").Select()
